$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 328, shifting existing rows 328:386 down to 329:387
$ws.Rows.Item(328).Insert()

# Populate the new row 328 with the new data record
$ws.Cells.Item(328, 1).Value2 = 4
$ws.Cells.Item(328, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(328, 3).Value = 'Los Lagos'
$ws.Cells.Item(328, 4).Value2 = 44995
$ws.Cells.Item(328, 5).Value2 = 10
$ws.Cells.Item(328, 6).Value = 'Fruta'
$ws.Cells.Item(328, 7).Value2 = 100108
$ws.Cells.Item(328, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(328, 9).Value2 = 100108005
$ws.Cells.Item(328, 10).Value = 'Piña'
$ws.Cells.Item(328, 11).Value = 'Caramelo'
$ws.Cells.Item(328, 12).Value = 'Primera'
$ws.Cells.Item(328, 13).Value2 = 300
$ws.Cells.Item(328, 14).Value2 = 24000
$ws.Cells.Item(328, 15).Value2 = 25000
$ws.Cells.Item(328, 16).Value2 = 24500
$ws.Cells.Item(328, 17).Value = '$/caja 12 unidades'
$ws.Cells.Item(328, 18).Value = 'Ecuador'
$ws.Cells.Item(328, 19).Value2 = 2042
$ws.Cells.Item(328, 20).Value2 = 12
